$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-43 all currently hold the date serial 45759
# (2025-04-12). Update them all to 45760 (2025-04-13).
$ws.Range("C2:C43").Value = 45760
